$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A3").Value = "TC002"
$ws.Range("B3").Value = "Samsung"
$ws.Range("C3").Value = "Tanmay"
$ws.Range("D3").Value = "Sarkar"

$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)

$ws.Range("D4").Select() | Out-Null
